# Normalize the composite-score group labels in column A (and their
# matching shared-string entries) to the lowercase/underscore style used
# by the other parameter groups (screen_time, sleep_time, ...).
#
#   Thinking          -> thinking
#   Body_signals      -> body_signals
#   Emotions          -> emotions
#   Composite scores  -> composite_score
#
# Using Cells.Replace (instead of writing .Value on each cell) edits the
# shared-string text in place, which is what happened in the real commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$null = $ws.Cells.Replace("Thinking", "thinking")
$null = $ws.Cells.Replace("Body_signals", "body_signals")
$null = $ws.Cells.Replace("Emotions", "emotions")
$null = $ws.Cells.Replace("Composite scores", "composite_score")

# Reselect the composite-score rows (A39:A41) as the last user action.
$null = $ws.Range("A39:A41").Select()
